# "D Tuan check Daotao" - add a new Đào tạo test-case row to the TestCaseAsign sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New test case row (row 13), mirrors the existing "Về TBD" test case (row 12)
# but for the "Đào tạo" page.
$ws.Range("A13").Value = "Check Đào tạo"
$ws.Range("B13").Value = "Kiểm tra link Về TDB"
$ws.Range("C13").Value = "WebElement: https://tbd.edu.vn/dao-tao/"
$ws.Range("D13").Value = "Test Link Passed!"
$ws.Range("E13").Value = "Đào tạo"
$ws.Range("F13").Value = "Pass"
$ws.Range("G13").Value = "Đình Tuấn"

# Column C (Test Data) now holds a long URL string - widen it to fit, same as
# Excel's own auto-fit behaviour when the new text was entered.
$ws.Columns.Item(3).ColumnWidth = 35.1667

# Leave the cursor where the author left off after entering the new row.
$ws.Range("F8").Select()
